# Rename the inline picture shapes' display names:
#   - BTec logo images (in the two headers): image1.jpg -> image2.jpg
#   - Pearson logo images (in the two footers): image2.png -> image1.png
#
# Both headers / footers (default + first-page) carry a copy of each logo,
# so every Section's Headers(1..2) / Footers(1..2) needs the same rename.

$d = $word.ActiveDocument

for ($secIdx = 1; $secIdx -le $d.Sections.Count; $secIdx++) {
    $section = $d.Sections.Item($secIdx)

    for ($hIdx = 1; $hIdx -le $section.Headers.Count; $hIdx++) {
        $header = $section.Headers.Item($hIdx)
        if ($header.Exists) {
            for ($shpIdx = 1; $shpIdx -le $header.Range.InlineShapes.Count; $shpIdx++) {
                $shape = $header.Range.InlineShapes.Item($shpIdx)
                if ($shape.AlternativeText -eq "BTec_Logo-Orange") {
                    $shape.Name = "image2.jpg"
                }
            }
        }
    }

    for ($fIdx = 1; $fIdx -le $section.Footers.Count; $fIdx++) {
        $footer = $section.Footers.Item($fIdx)
        if ($footer.Exists) {
            for ($shpIdx = 1; $shpIdx -le $footer.Range.InlineShapes.Count; $shpIdx++) {
                $shape = $footer.Range.InlineShapes.Item($shpIdx)
                if ($shape.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shape.Name = "image1.png"
                }
            }
        }
    }
}
